$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "41.797.84"
$ws.Cells.Item(2, 5).Value = "  +2.21%  "
$ws.Cells.Item(3, 4).Value = "2.228.89"
$ws.Cells.Item(3, 5).Value = "  +0.10%  "
$ws.Cells.Item(4, 5).Value = "  +0.00%  "
$ws.Cells.Item(5, 4).Value = "'232.08"
$ws.Cells.Item(5, 5).Value = "  +1.25%  "
$ws.Cells.Item(6, 4).Value = "'0.624"
$ws.Cells.Item(6, 5).Value = "  -1.61%  "
$ws.Cells.Item(7, 4).Value = "'60.49"
$ws.Cells.Item(7, 5).Value = "  -7.00%  "
$ws.Cells.Item(8, 5).Value = "  -0.01%  "
$ws.Cells.Item(9, 5).Value = "  -0.65%  "
$ws.Cells.Item(10, 4).Value = "'58.19"
$ws.Cells.Item(10, 5).Value = "  -1.91%  "
$ws.Cells.Item(11, 4).Value = "'0.0900"
$ws.Cells.Item(11, 5).Value = "  +2.12%  "
$ws.Cells.Item(12, 5).Value = "  -0.34%  "
$ws.Cells.Item(13, 4).Value = "2.561.29"
$ws.Cells.Item(13, 5).Value = "  +0.37%  "
$ws.Cells.Item(14, 4).Value = "'15.60"
$ws.Cells.Item(14, 5).Value = "  -3.39%  "
$ws.Cells.Item(15, 4).Value = "'22.70"
$ws.Cells.Item(15, 5).Value = "  +1.13%  "
$ws.Cells.Item(16, 5).Value = "  -2.87%  "
$ws.Cells.Item(17, 4).Value = "'5.62"
$ws.Cells.Item(17, 5).Value = "  -0.51%  "
$ws.Cells.Item(18, 4).Value = "2.241.07"
$ws.Cells.Item(18, 5).Value = "  +0.91%  "
$ws.Cells.Item(19, 4).Value = "41.738.77"
$ws.Cells.Item(19, 5).Value = "  +2.29%  "
$ws.Cells.Item(20, 4).Value = "0.0₃0911"
$ws.Cells.Item(20, 5).Value = "  +0.30%  "
$ws.Cells.Item(21, 4).Value = "'72.52"
$ws.Cells.Item(21, 5).Value = "  -2.11%  "
$ws.Cells.Item(22, 4).Value = "'6.16"
$ws.Cells.Item(22, 5).Value = "  -0.14%  "
$ws.Cells.Item(23, 4).Value = "'248.13"
$ws.Cells.Item(23, 5).Value = "  -2.92%  "
$ws.Cells.Item(24, 5).Value = "  -0.16%  "
$ws.Cells.Item(25, 2).Value = "Toncoin"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(25, 4).Value = "'2.38"
$ws.Cells.Item(25, 5).Value = "  +2.88%  "
$ws.Cells.Item(26, 2).Value = "PancakeSwap"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(26, 4).Value = "'2.39"
$ws.Cells.Item(26, 5).Value = "  -0.09%  "
$ws.Cells.Item(27, 4).Value = "'9.61"
$ws.Cells.Item(27, 5).Value = "  -1.55%  "
$ws.Cells.Item(28, 4).Value = "'169.61"
$ws.Cells.Item(28, 5).Value = "  -2.12%  "
$ws.Cells.Item(29, 5).Value = "  -2.09%  "
$ws.Cells.Item(30, 4).Value = "'19.93"
$ws.Cells.Item(30, 5).Value = "  -2.28%  "
$ws.Cells.Item(31, 5).Value = "  -2.93%  "
$ws.Cells.Item(32, 5).Value = "  -8.06%  "
$ws.Cells.Item(33, 5).Value = "  -1.42%  "
$ws.Cells.Item(34, 4).Value = "'5.02"
$ws.Cells.Item(34, 5).Value = "  +3.52%  "
$ws.Cells.Item(35, 4).Value = "'4.70"
$ws.Cells.Item(35, 5).Value = "  +0.43%  "
$ws.Cells.Item(36, 4).Value = "'0.0654"
$ws.Cells.Item(36, 5).Value = "  +3.15%  "
$ws.Cells.Item(37, 5).Value = "  -9.22%  "
$ws.Cells.Item(38, 5).Value = "  -5.31%  "
$ws.Cells.Item(39, 5).Value = "  -4.17%  "
$ws.Cells.Item(40, 5).Value = "  +0.02%  "
$ws.Cells.Item(41, 2).Value = "TerraClassic"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Cells.Item(41, 4).Value = "'0.000236"
$ws.Cells.Item(41, 5).Value = "  +3.08%  "
$ws.Cells.Item(42, 2).Value = "VeChain"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(42, 4).Value = "'0.0240"
$ws.Cells.Item(42, 5).Value = "  +1.21%  "
$ws.Cells.Item(43, 4).Value = "'8.60"
$ws.Cells.Item(43, 5).Value = "  -1.80%  "
$ws.Cells.Item(44, 4).Value = "'1.23"
$ws.Cells.Item(44, 5).Value = "  -1.78%  "
$ws.Cells.Item(45, 4).Value = "'98.78"
$ws.Cells.Item(45, 5).Value = "  -3.44%  "
$ws.Cells.Item(46, 4).Value = "'4.46"
$ws.Cells.Item(46, 5).Value = "  -9.79%  "
$ws.Cells.Item(47, 4).Value = "'0.0964"
$ws.Cells.Item(47, 5).Value = "  +1.80%  "
$ws.Cells.Item(48, 4).Value = "1.471.24"
$ws.Cells.Item(48, 5).Value = "  -3.02%  "
$ws.Cells.Item(49, 4).Value = "'16.60"
$ws.Cells.Item(49, 5).Value = "  -6.01%  "
$ws.Cells.Item(50, 4).Value = "'2.31"
$ws.Cells.Item(50, 5).Value = "  +9.05%  "
$ws.Cells.Item(51, 4).Value = "'2.79"
$ws.Cells.Item(51, 5).Value = "  -1.48%  "
